$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1914893617021277
$ws.Range("C2").Value = 0.5585106382978723
$ws.Range("J2").Value = 0.02127659574468085
$ws.Range("P2").Value = 0.1542553191489362
$ws.Range("S2").Value = 0.07446808510638298
$ws.Range("B3").Value = 0.01401869158878505
$ws.Range("C3").Value = 0.01869158878504673
$ws.Range("J3").Value = 0.02336448598130841
$ws.Range("O3").Value = 0.004672897196261682
$ws.Range("P3").Value = 0.719626168224299
$ws.Range("S3").Value = 0.2196261682242991
$ws.Range("J4").Value = 0.1020408163265306
$ws.Range("P4").Value = 0.5918367346938775
$ws.Range("S4").Value = 0.3061224489795918
$ws.Range("B6").Value = 0.07623318385650224
$ws.Range("D6").Value = 0.02242152466367713
$ws.Range("F6").Value = 0.07623318385650224
$ws.Range("J6").Value = 0.2466367713004484
$ws.Range("O6").Value = 0.03139013452914798
$ws.Range("Q6").Value = 0.08071748878923767
$ws.Range("R6").Value = 0.1165919282511211
$ws.Range("S6").Value = 0.3497757847533632
$ws.Range("B7").Value = 0.1434262948207171
$ws.Range("D7").Value = 0.0199203187250996
$ws.Range("F7").Value = 0.05976095617529881
$ws.Range("J7").Value = 0.1713147410358566
$ws.Range("O7").Value = 0.03187250996015936
$ws.Range("Q7").Value = 0.1434262948207171
$ws.Range("R7").Value = 0.1035856573705179
$ws.Range("S7").Value = 0.3266932270916335
$ws.Range("B8").Value = 0.08108108108108109
$ws.Range("D8").Value = 0.01351351351351351
$ws.Range("E8").Value = 0.001930501930501931
$ws.Range("F8").Value = 0.05019305019305019
$ws.Range("J8").Value = 0.1370656370656371
$ws.Range("O8").Value = 0.01737451737451737
$ws.Range("Q8").Value = 0.1621621621621622
$ws.Range("R8").Value = 0.1081081081081081
$ws.Range("S8").Value = 0.4285714285714285
$ws.Range("B9").Value = 0.1146245059288538
$ws.Range("D9").Value = 0.01976284584980237
$ws.Range("F9").Value = 0.05928853754940711
$ws.Range("J9").Value = 0.1106719367588933
$ws.Range("O9").Value = 0.003952569169960474
$ws.Range("Q9").Value = 0.1778656126482213
$ws.Range("R9").Value = 0.1027667984189723
$ws.Range("S9").Value = 0.4110671936758893
$ws.Range("B10").Value = 0.1199451679232351
$ws.Range("D10").Value = 0.0205620287868403
$ws.Range("E10").Value = 0.0006854009595613434
$ws.Range("F10").Value = 0.05962988348183688
$ws.Range("J10").Value = 0.1425633995887594
$ws.Range("O10").Value = 0.01713502398903358
$ws.Range("Q10").Value = 0.1706648389307745
$ws.Range("R10").Value = 0.1062371487320082
$ws.Range("S10").Value = 0.3625771076079506
$ws.Range("G11").Value = 0.1906005221932115
$ws.Range("J11").Value = 0.07571801566579635
$ws.Range("K11").Value = 0.2036553524804177
$ws.Range("L11").Value = 0.5169712793733682
$ws.Range("S11").Value = 0.01305483028720627
$ws.Range("G12").Value = 0.7205882352941176
$ws.Range("J12").Value = 0.1764705882352941
$ws.Range("K12").Value = 0.009803921568627451
$ws.Range("L12").Value = 0.02941176470588235
$ws.Range("S12").Value = 0.06372549019607843
$ws.Range("F13").Value = 0.01639344262295082
$ws.Range("G13").Value = 0.6721311475409836
$ws.Range("J13").Value = 0.2459016393442623
$ws.Range("S13").Value = 0.06557377049180328
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.1244635193133047
$ws.Range("I15").Value = 0.06437768240343347
$ws.Range("J15").Value = 0.3433476394849785
$ws.Range("K15").Value = 0.09871244635193133
$ws.Range("M15").Value = 0.02145922746781116
$ws.Range("O15").Value = 0.04721030042918455
$ws.Range("S15").Value = 0.2832618025751073
$ws.Range("F16").Value = 0.01271186440677966
$ws.Range("H16").Value = 0.1822033898305085
$ws.Range("I16").Value = 0.09745762711864407
$ws.Range("J16").Value = 0.3813559322033898
$ws.Range("K16").Value = 0.1186440677966102
$ws.Range("M16").Value = 0.02542372881355932
$ws.Range("O16").Value = 0.06779661016949153
$ws.Range("S16").Value = 0.1144067796610169
$ws.Range("F17").Value = 0.01834862385321101
$ws.Range("H17").Value = 0.1903669724770642
$ws.Range("I17").Value = 0.1100917431192661
$ws.Range("J17").Value = 0.3692660550458716
$ws.Range("K17").Value = 0.1422018348623853
$ws.Range("M17").Value = 0.02752293577981652
$ws.Range("O17").Value = 0.03669724770642202
$ws.Range("S17").Value = 0.1055045871559633
$ws.Range("F18").Value = 0.02097902097902098
$ws.Range("H18").Value = 0.1748251748251748
$ws.Range("I18").Value = 0.1048951048951049
$ws.Range("J18").Value = 0.4090909090909091
$ws.Range("K18").Value = 0.0944055944055944
$ws.Range("M18").Value = 0.02097902097902098
$ws.Range("O18").Value = 0.05594405594405594
$ws.Range("S18").Value = 0.1188811188811189
$ws.Range("F19").Value = 0.008241758241758242
$ws.Range("H19").Value = 0.2190934065934066
$ws.Range("I19").Value = 0.09409340659340659
$ws.Range("J19").Value = 0.3592032967032967
$ws.Range("K19").Value = 0.1085164835164835
$ws.Range("M19").Value = 0.02335164835164835
$ws.Range("N19").Value = 0.00206043956043956
$ws.Range("O19").Value = 0.06387362637362637
$ws.Range("S19").Value = 0.1215659340659341
